# Natmi following Dr Hou advice
#
# Re-running the NATMI ligand-receptor pairing for Artn-Gfra1 after adding
# the "ECs" sending/target cluster to the analysis. Rows 2-4 (the existing
# FAPs/sCs combinations) get their expression/specificity statistics
# recomputed against the larger cluster set, and rows 5-7 are appended for
# the new combinations that involve "ECs".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns A:T per the header row:
# Sending cluster, Ligand symbol, Receptor symbol, Target cluster,
# Ligand-expressing cells, Ligand detection rate,
# Ligand average expression value, Ligand total expression value,
# Ligand derived specificity of average expression value,
# Ligand derived specificity of total expression value,
# Receptor-expressing cells, Receptor detection rate,
# Receptor average expression value, Receptor total expression value,
# Receptor derived specificity of average expression value,
# Receptor derived specificity of total expression value,
# Edge average expression weight, Edge total expression weight,
# Edge average expression derived specificity, Edge total expression derived specificity
$data = @{
    2 = @("FAPs", "Artn", "Gfra1", "ECs", 3, 1, 2.306751666666667, 6.920255, 0.8617934705859058, 0.8617934705859057, 2, 0.6666666666666666, 0.3392376666666667, 1.017713, 0.01394884456468018, 0.01394884456468017, 0.7825370529794445, 7.042833476815001, 0.01202102316805908, 0.01202102316805907)
    3 = @("FAPs", "Artn", "Gfra1", "FAPs", 3, 1, 2.306751666666667, 6.920255, 0.8617934705859058, 0.8617934705859057, 3, 1, 20.328499, 60.985497, 0.8358714277529804, 0.8358714277529804, 46.89279894908166, 422.0351905417349, 0.7203485386868372, 0.7203485386868371)
    4 = @("FAPs", "Artn", "Gfra1", "sCs", 3, 1, 2.306751666666667, 6.920255, 0.8617934705859058, 0.8617934705859057, 3, 1, 3.65239, 10.95717, 0.1501797276823394, 0.1501797276823394, 8.425156719816666, 75.82641047835, 0.1294239087310095, 0.1294239087310095)
    5 = @("sCs", "Artn", "Gfra1", "ECs", 3, 1, 0.3699356666666667, 1.109807, 0.1382065294140942, 0.1382065294140942, 2, 0.6666666666666666, 0.3392376666666667, 1.017713, 0.01394884456468018, 0.01394884456468017, 0.1254961123767778, 1.129465011391, 0.001927821396621099, 0.001927821396621098)
    6 = @("sCs", "Artn", "Gfra1", "FAPs", 3, 1, 0.3699356666666667, 1.109807, 0.1382065294140942, 0.1382065294140942, 3, 1, 20.328499, 60.985497, 0.8358714277529804, 0.8358714277529804, 7.520236829897666, 67.682131469079, 0.1155228890661432, 0.1155228890661432)
    7 = @("sCs", "Artn", "Gfra1", "sCs", 3, 1, 0.3699356666666667, 1.109807, 0.1382065294140942, 0.1382065294140942, 3, 1, 3.65239, 10.95717, 0.1501797276823394, 0.1501797276823394, 1.351149329576667, 12.16034396619, 0.0207558189513299, 0.02075581895132989)
}

foreach ($r in 2..7) {
    $rowValues = $data[$r]
    for ($col = 1; $col -le $rowValues.Length; $col++) {
        $ws.Cells.Item($r, $col).Value = $rowValues[$col - 1]
    }
}
